$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header cell in H1, copying the existing header formatting
# (bold font + border + centered alignment) from the neighboring G1 cell.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Fill in the Save column values for rows 2-7 (1 = saved, 0 = not saved)
$saveValues = @(1, 1, 0, 1, 0, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
